$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header cell F1 with value then copy formatting (style) from E1
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Populate F2:F95 with time_taken timestamps
$ws.Range("F2").Value = "2021-10-05 10:52:32.967637"
$ws.Range("F3").Value = "2021-10-05 10:52:32.967648"
$ws.Range("F4").Value = "2021-10-05 10:52:32.967652"
$ws.Range("F5").Value = "2021-10-05 10:52:32.967654"
$ws.Range("F6").Value = "2021-10-05 10:52:32.967657"
$ws.Range("F7").Value = "2021-10-05 10:52:32.967660"
$ws.Range("F8").Value = "2021-10-05 10:52:32.967662"
$ws.Range("F9").Value = "2021-10-05 10:52:32.967665"
$ws.Range("F10").Value = "2021-10-05 10:52:32.967667"
$ws.Range("F11").Value = "2021-10-05 10:52:32.967670"
$ws.Range("F12").Value = "2021-10-05 10:52:32.967672"
$ws.Range("F13").Value = "2021-10-05 10:52:32.967675"
$ws.Range("F14").Value = "2021-10-05 10:52:32.967677"
$ws.Range("F15").Value = "2021-10-05 10:52:32.967680"
$ws.Range("F16").Value = "2021-10-05 10:52:32.967682"
$ws.Range("F17").Value = "2021-10-05 10:52:32.967685"
$ws.Range("F18").Value = "2021-10-05 10:52:32.967688"
$ws.Range("F19").Value = "2021-10-05 10:52:32.967690"
$ws.Range("F20").Value = "2021-10-05 10:52:32.967693"
$ws.Range("F21").Value = "2021-10-05 10:52:32.967695"
$ws.Range("F22").Value = "2021-10-05 10:52:32.967698"
$ws.Range("F23").Value = "2021-10-05 10:52:32.967700"
$ws.Range("F24").Value = "2021-10-05 10:52:32.967703"
$ws.Range("F25").Value = "2021-10-05 10:52:32.967705"
$ws.Range("F26").Value = "2021-10-05 10:52:32.967708"
$ws.Range("F27").Value = "2021-10-05 10:52:32.967711"
$ws.Range("F28").Value = "2021-10-05 10:52:32.967713"
$ws.Range("F29").Value = "2021-10-05 10:52:32.967716"
$ws.Range("F30").Value = "2021-10-05 10:52:32.967718"
$ws.Range("F31").Value = "2021-10-05 10:52:32.967720"
$ws.Range("F32").Value = "2021-10-05 10:52:32.967723"
$ws.Range("F33").Value = "2021-10-05 10:52:32.967725"
$ws.Range("F34").Value = "2021-10-05 10:52:32.967729"
$ws.Range("F35").Value = "2021-10-05 10:52:32.967731"
$ws.Range("F36").Value = "2021-10-05 10:52:32.967734"
$ws.Range("F37").Value = "2021-10-05 10:52:32.967736"
$ws.Range("F38").Value = "2021-10-05 10:52:32.967739"
$ws.Range("F39").Value = "2021-10-05 10:52:32.967741"
$ws.Range("F40").Value = "2021-10-05 10:52:32.967743"
$ws.Range("F41").Value = "2021-10-05 10:52:32.967746"
$ws.Range("F42").Value = "2021-10-05 10:52:32.967749"
$ws.Range("F43").Value = "2021-10-05 10:52:32.967751"
$ws.Range("F44").Value = "2021-10-05 10:52:32.967754"
$ws.Range("F45").Value = "2021-10-05 10:52:32.967756"
$ws.Range("F46").Value = "2021-10-05 10:52:32.967759"
$ws.Range("F47").Value = "2021-10-05 10:52:32.967761"
$ws.Range("F48").Value = "2021-10-05 10:52:32.967763"
$ws.Range("F49").Value = "2021-10-05 10:52:32.967766"
$ws.Range("F50").Value = "2021-10-05 10:52:32.967768"
$ws.Range("F51").Value = "2021-10-05 10:52:32.967771"
$ws.Range("F52").Value = "2021-10-05 10:52:32.967773"
$ws.Range("F53").Value = "2021-10-05 10:52:32.967776"
$ws.Range("F54").Value = "2021-10-05 10:52:32.967779"
$ws.Range("F55").Value = "2021-10-05 10:52:32.967781"
$ws.Range("F56").Value = "2021-10-05 10:52:32.967783"
$ws.Range("F57").Value = "2021-10-05 10:52:32.967786"
$ws.Range("F58").Value = "2021-10-05 10:52:32.967789"
$ws.Range("F59").Value = "2021-10-05 10:52:32.967791"
$ws.Range("F60").Value = "2021-10-05 10:52:32.967794"
$ws.Range("F61").Value = "2021-10-05 10:52:32.967796"
$ws.Range("F62").Value = "2021-10-05 10:52:32.967799"
$ws.Range("F63").Value = "2021-10-05 10:52:32.967801"
$ws.Range("F64").Value = "2021-10-05 10:52:32.967803"
$ws.Range("F65").Value = "2021-10-05 10:52:32.967806"
$ws.Range("F66").Value = "2021-10-05 10:52:32.967809"
$ws.Range("F67").Value = "2021-10-05 10:52:32.967812"
$ws.Range("F68").Value = "2021-10-05 10:52:32.967815"
$ws.Range("F69").Value = "2021-10-05 10:52:32.967817"
$ws.Range("F70").Value = "2021-10-05 10:52:32.967820"
$ws.Range("F71").Value = "2021-10-05 10:52:32.967822"
$ws.Range("F72").Value = "2021-10-05 10:52:32.967825"
$ws.Range("F73").Value = "2021-10-05 10:52:32.967827"
$ws.Range("F74").Value = "2021-10-05 10:52:32.967830"
$ws.Range("F75").Value = "2021-10-05 10:52:32.967832"
$ws.Range("F76").Value = "2021-10-05 10:52:32.967835"
$ws.Range("F77").Value = "2021-10-05 10:52:32.967837"
$ws.Range("F78").Value = "2021-10-05 10:52:32.967842"
$ws.Range("F79").Value = "2021-10-05 10:52:32.967845"
$ws.Range("F80").Value = "2021-10-05 10:52:32.967847"
$ws.Range("F81").Value = "2021-10-05 10:52:32.967850"
$ws.Range("F82").Value = "2021-10-05 10:52:32.967852"
$ws.Range("F83").Value = "2021-10-05 10:52:32.967855"
$ws.Range("F84").Value = "2021-10-05 10:52:32.967857"
$ws.Range("F85").Value = "2021-10-05 10:52:32.967860"
$ws.Range("F86").Value = "2021-10-05 10:52:32.967862"
$ws.Range("F87").Value = "2021-10-05 10:52:32.967865"
$ws.Range("F88").Value = "2021-10-05 10:52:32.967867"
$ws.Range("F89").Value = "2021-10-05 10:52:32.967870"
$ws.Range("F90").Value = "2021-10-05 10:52:32.967872"
$ws.Range("F91").Value = "2021-10-05 10:52:32.967875"
$ws.Range("F92").Value = "2021-10-05 10:52:32.967877"
$ws.Range("F93").Value = "2021-10-05 10:52:32.967880"
$ws.Range("F94").Value = "2021-10-05 10:52:32.967884"
$ws.Range("F95").Value = "2021-10-05 10:52:32.967887"

$excel.CutCopyMode = 0
